$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work_Sheet")

# ---------------------------------------------------------------------------
# Row 14: E14 already holds 1 but with no explicit style; the target gives it
# the column's standard "General" style (s="8"). Clearing then re-writing the
# value forces the cell to be re-created and pick up the column default style.
# ---------------------------------------------------------------------------
$ws.Range("E14").ClearContents()
$ws.Range("E14").Value = 1

# ---------------------------------------------------------------------------
# Row 15: only the date changes (2019-10-16 -> 2019-10-18).
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 43756

# ---------------------------------------------------------------------------
# Row 27: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 43756
$ws.Range("E27").Value = 1

# ---------------------------------------------------------------------------
# Row 43: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C43").Value = 43756
$ws.Range("E43").Value = 1

# ---------------------------------------------------------------------------
# Row 44: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C44").Value = 43756
$ws.Range("E44").Value = 1

# ---------------------------------------------------------------------------
# Row 48: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C48").Value = 43756
$ws.Range("E48").Value = 1

# ---------------------------------------------------------------------------
# Row 49: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C49").Value = 43756
$ws.Range("E49").Value = 1

# ---------------------------------------------------------------------------
# Row 56: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C56").Value = 43756
$ws.Range("E56").Value = 1

# ---------------------------------------------------------------------------
# Row 59: only the date changes (2019-10-16 -> 2019-10-18).
# ---------------------------------------------------------------------------
$ws.Range("C59").Value = 43756

# ---------------------------------------------------------------------------
# Row 60: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C60").Value = 43756
$ws.Range("E60").Value = 1

# ---------------------------------------------------------------------------
# Row 63: date changes, and the existing Days=1 cell gains the standard style.
# ---------------------------------------------------------------------------
$ws.Range("C63").Value = 43756
$ws.Range("E63").ClearContents()
$ws.Range("E63").Value = 1

# ---------------------------------------------------------------------------
# Row 64: date stays the same; the existing Days=1 cell gains the standard
# style.
# ---------------------------------------------------------------------------
$ws.Range("E64").ClearContents()
$ws.Range("E64").Value = 1

# ---------------------------------------------------------------------------
# Row 70: only the date changes (2019-10-16 -> 2019-10-18).
# ---------------------------------------------------------------------------
$ws.Range("C70").Value = 43756

# ---------------------------------------------------------------------------
# Row 71: only the date changes (2019-10-16 -> 2019-10-18).
# ---------------------------------------------------------------------------
$ws.Range("C71").Value = 43756

# ---------------------------------------------------------------------------
# Row 83: date gets filled in, and a new Days=1 cell is added.
# ---------------------------------------------------------------------------
$ws.Range("C83").Value = 43756
$ws.Range("E83").Value = 1

# ---------------------------------------------------------------------------
# Row 87: date stays the same; the existing Days=1 cell gains the standard
# style.
# ---------------------------------------------------------------------------
$ws.Range("E87").ClearContents()
$ws.Range("E87").Value = 1

# ---------------------------------------------------------------------------
# Row 89: date stays the same; the existing Days=1 cell gains the standard
# style.
# ---------------------------------------------------------------------------
$ws.Range("E89").ClearContents()
$ws.Range("E89").Value = 1

# ---------------------------------------------------------------------------
# Row 94: date changes; the existing Days=2 cell gains the standard style.
# ---------------------------------------------------------------------------
$ws.Range("C94").Value = 43756
$ws.Range("E94").ClearContents()
$ws.Range("E94").Value = 2
